# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @{
    2  = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    3  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    4  = @(0.04172184405617529, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1.60109356927828)
    5  = @(0.2881169905109251, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.598097515653722)
    6  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    7  = @(1.445647641019636, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 3.005019366241741)
    8  = @(0.003078177322033415, 0.002658071450198252, 3.223369029078222, 0.5333859586016987, 3.762491236452152)
    9  = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 24.14949828602258)
    10 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    11 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    12 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    13 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 18.91276827552123)
    14 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    15 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    16 = @(0.6545652718822623, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 20.20963986881668)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
